$wb = $excel.ActiveWorkbook

# --- "tab" sheet: move the selection off J18 onto A4 and deselect the tab ---
# (do this before creating/activating the new sheet so "tab" loses tabSelected)
$tabWs = $wb.Worksheets.Item("tab")
$tabWs.Range("A4").Select() | Out-Null

# --- clone "tab" into a new worksheet placed right after it ---
$tabWs.Copy($null, $tabWs) | Out-Null
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "var_set"

# Insert a fresh row at 3 (pushes the "url"/"title" rows down one, taking their
# formatting with them) to make room for the new "try 1" step row.
$newWs.Rows.Item(3).Insert() | Out-Null

# Drop the old trailing "tab:back"/"title?" rows (12-13 after the insert shift);
# the var_set case only needs 11 data rows.
$newWs.Range("A12:E13").EntireRow.Delete() | Out-Null

# --- rewrite the step data for the var action test case ---
$newWs.Range("A2").Value = "sheet - action - var"
$newWs.Range("B3").Value = "try 1"

$newWs.Range("C6").Value = "button#multi"
$newWs.Range("D6").Value = "var"
$newWs.Range("E6").Value = "myvar1"

$newWs.Range("D7").Value = "print"
$newWs.Range("E7").Value = "something"

$newWs.Range("C8").ClearContents() | Out-Null
$newWs.Range("D8").Value = "print"
$newWs.Range("E8").Value = "xx{{myvar1}}"

$newWs.Range("C9").Value = "multi"
$newWs.Range("D9").Value = "var:set"
$newWs.Range("E9").Value = "myvar2"

$newWs.Range("D10").Value = "print"
$newWs.Range("E10").Value = "{{myvar2}}"

$newWs.Range("C11").Value = "button#{{myvar2}}"
$newWs.Range("D11").Value = "assert"
$newWs.Range("E11").Value = "Muiltiple windows"

# keep a trailing blank row marker (mirrors the other step sheets) and put the
# cursor where the authored workbook leaves it
$newWs.Rows.Item(24).RowHeight = 13.8
$newWs.Range("A3").Select() | Out-Null
